$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).ColumnWidth = 15.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.833333333333334

$ws.Range("A1").Value = -0.072935418983419709
$ws.Range("B1").Value = 0.072813847406081322
$ws.Range("A2").Value = -0.00063708338208279258
$ws.Range("B2").Value = 0.00014888658588319004
$ws.Range("A3").Value = 0.10278360103213302
$ws.Range("B3").Value = -0.10315800843015666
$ws.Range("A4").Value = -0.18083865424986101
$ws.Range("B4").Value = 0.17995019019635095
$ws.Range("A5").Value = -0.17395019038274384
$ws.Range("B5").Value = 0.17216960330022335
$ws.Range("A6").Value = -0.059879446679754977
$ws.Range("B6").Value = 0.059834267735672242
$ws.Range("A7").Value = -0.039834267962124414
$ws.Range("B7").Value = 0.039761272083040922
$ws.Range("A8").Value = -0.019761272310697464
$ws.Range("B8").Value = 0.019730770887114168
$ws.Range("A9").Value = -0.061488524460273197
$ws.Range("B9").Value = 0.061178403399488879
$ws.Range("A10").Value = -0.055178403597807346
$ws.Range("B10").Value = 0.055132332857048993
$ws.Range("A11").Value = -0.051473598501129203
$ws.Range("B11").Value = 0.051394321058346293
$ws.Range("A12").Value = -0.045394321257450354
$ws.Range("B12").Value = 0.045151483778357981
$ws.Range("A13").Value = -0.039151483980573332
$ws.Range("B13").Value = 0.039085215650225713
$ws.Range("A14").Value = -0.027085215867558077
$ws.Range("B14").Value = 0.027052984899937904
$ws.Range("A15").Value = -0.021052985103785282
$ws.Range("B15").Value = 0.021027596391659742
$ws.Range("A16").Value = -0.015027596596164372
$ws.Range("B16").Value = 0.015004643515456717
$ws.Range("A17").Value = -0.0090046437208419761
$ws.Range("B17").Value = 0.0089999997873313475
$ws.Range("A18").Value = -0.062176764218133229
$ws.Range("B18").Value = 0.062138303200473644
$ws.Range("A19").Value = -0.053138303392644026
$ws.Range("B19").Value = 0.05286701199510313
$ws.Range("A20").Value = -0.018014084048367351
$ws.Range("B20").Value = 0.018004314040345193
$ws.Range("A21").Value = -0.0090043142362441486
$ws.Range("B21").Value = 0.0089999998039242968
$ws.Range("A22").Value = -0.093932487810018728
$ws.Range("B22").Value = 0.093624102184801217
$ws.Range("A23").Value = -0.084624102378434429
$ws.Range("B23").Value = 0.084124846493379479
$ws.Range("A24").Value = -0.042124846768051505
$ws.Range("B24").Value = 0.041999999723906889
$ws.Range("A25").Value = -0.09333737643110851
$ws.Range("B25").Value = 0.093166431159701801
$ws.Range("A26").Value = -0.087166431353214335
$ws.Range("B26").Value = 0.086951867392862425
$ws.Range("A27").Value = -0.080951867587313764
$ws.Range("B27").Value = 0.080239305606954137
$ws.Range("A28").Value = -0.074239305805482658
$ws.Range("B28").Value = 0.073760615066041701
$ws.Range("A29").Value = -0.061760615281293951
$ws.Range("B29").Value = 0.061540814051348036
$ws.Range("A30").Value = -0.04216708291430793
$ws.Range("B30").Value = 0.042020129415034635
$ws.Range("A31").Value = -0.027020129640584756
$ws.Range("B31").Value = 0.027000929783906003
$ws.Range("A32").Value = -0.0060009300239682872
$ws.Range("B32").Value = 0.0059999997953248396
